# updated NA US scripts
# Renames the "-test" suffixed hostess/cohostess/guest values to "-test1"
# on the "cred" sheet of virtualShowData.xlsx, and moves the active
# selection/scroll position from column K to column J/G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds the "test" data used by the NA US Katalon scripts.
# Update each "-testX abc" / "abc, Xtest" value to use the "test1" suffix.
$ws.Range("B2").Value = "ahostess-test1 abc"
$ws.Range("C2").Value = "bcohost-test1 abc"
$ws.Range("D2").Value = "guest1-test1 abc"
$ws.Range("E2").Value = "guest2-test1 abc"
$ws.Range("F2").Value = "abc, ahostess-test1"
$ws.Range("G2").Value = "abc, bcohost-test1"
$ws.Range("I2").Value = "abc, guest1-test1"
$ws.Range("J2").Value = "abc, guest2-test1"

# Update the saved view/selection state: scroll to column G and select J2.
$ws.Range("J2").Select()
